# Add a new worksheet "Uni_SPI1" that holds a plain-value copy of the
# "SPI1" column (column B, header + 150 observations) from sheet "SPI1".

$wb = $excel.ActiveWorkbook

# "SPI1" is the second sheet (SPI3, SPI1) in the workbook.
$spi1 = $wb.Worksheets.Item(2)

# Select & copy the SPI1 column (B1:B151 contains the header + data;
# matches the author selecting column B before copying it out).
$spi1.Activate()
$spi1.Range("B1:B151").Copy()

# New sheet goes after the last existing sheet, so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Uni_SPI1"

# Paste as values only into column A of the new sheet (no source styles).
$newSheet.Range("A1").PasteSpecial(-4163)

# Restore SPI1's on-sheet selection to the whole of column B, reflecting
# the copy operation that sourced the new sheet's data.
$spi1.Activate()
$spi1.Columns("B").Select()

# Leave the new sheet active, with its last-used selection.
$newSheet.Activate()
$newSheet.Range("G18").Select()
